$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 170599.5017100554
$ws.Range("C2").Value = 243633.8521784183
$ws.Range("D2").Value = 288073.0325529571
$ws.Range("E2").Value = 320396.3342176132
$ws.Range("B3").Value = 210101.367450521
$ws.Range("C3").Value = 298388.0969814305
$ws.Range("D3").Value = 346302.7183486127
$ws.Range("E3").Value = 385064.3886487246
$ws.Range("B4").Value = 175774.1323965691
$ws.Range("C4").Value = 255451.2416646481
$ws.Range("D4").Value = 304989.6595546202
$ws.Range("E4").Value = 346963.6986408023
$ws.Range("B5").Value = 147650.4968770791
$ws.Range("C5").Value = 207747.4917195718
$ws.Range("D5").Value = 234416.135689859
$ws.Range("E5").Value = 259393.1022128583
$ws.Range("B6").Value = 128974.6871793733
$ws.Range("C6").Value = 180780.1452833392
$ws.Range("D6").Value = 205760.4743610043
$ws.Range("E6").Value = 225446.1385113747
$ws.Range("B7").Value = 14094.52972051891
$ws.Range("C7").Value = 19591.65091593164
$ws.Range("D7").Value = 22351.40199668508
$ws.Range("E7").Value = 24216.33805023737
$ws.Range("B8").Value = 685732.7234247532
$ws.Range("C8").Value = 978188.8231966568
$ws.Range("D8").Value = 1149251.332961092
$ws.Range("E8").Value = 1254841.582206369
$ws.Range("B9").Value = 195037.3496874791
$ws.Range("C9").Value = 272922.2820446609
$ws.Range("D9").Value = 309747.1584284796
$ws.Range("E9").Value = 337036.9602536431
$ws.Range("B10").Value = 83730.22944576826
$ws.Range("C10").Value = 113392.3348628892
$ws.Range("D10").Value = 129405.5374664847
$ws.Range("E10").Value = 137288.767285568
$ws.Range("B11").Value = 15259.14189030051
$ws.Range("C11").Value = 19581.7296019696
$ws.Range("D11").Value = 22083.70857085698
$ws.Range("E11").Value = 25295.97051138061
$ws.Range("B12").Value = 34229.96557708149
$ws.Range("C12").Value = 46490.64867722927
$ws.Range("D12").Value = 50729.55331958831
$ws.Range("E12").Value = 52421.17798267669
$ws.Range("B13").Value = 46934.74596435264
$ws.Range("C13").Value = 62908.47048988931
$ws.Range("D13").Value = 72328.30879069919
$ws.Range("E13").Value = 77603.82159470905
